$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'Energy and Environment Impacts'
    'Pollution'
    'Environmental Science'
    'Physical Sciences'
    'Employment and Welfare Studies'
    'General Health Professions'
    'Health Professions'
    'Health Sciences'
    'Agriculture and Farm Safety'
    'Plant Science'
    'Agricultural and Biological Sciences'
    'Life Sciences'
    'Mental Health Treatment and Access'
    'Social Psychology'
    'Psychology'
    'Social Sciences'
    'Cultural Differences and Values'
    'Community Health and Development'
    'Diabetes Management and Education'
    'Endocrinology, Diabetes and Metabolism'
    'Medicine'
    'Chronic Disease Management Strategies'
    'Epidemiology'
    'Cardiac Health and Mental Health'
    'Cardiology and Cardiovascular Medicine'
    'Participatory Visual Research Methods'
    'Sociology and Political Science'
    'Air Quality and Health Impacts'
    'Health, Toxicology and Mutagenesis'
    'Menopause: Health Impacts and Treatments'
    'Health and Wellbeing Research'
    'Sleep and Work-Related Fatigue'
    'Experimental and Cognitive Psychology'
    'Cervical Cancer and HPV Research'
    'Vaccine Coverage and Hesitancy'
    'Health'
    'Hepatitis B Virus Studies'
    'Child and Adolescent Psychosocial and Emotional Development'
    'Clinical Psychology'
    'Maternal Mental Health During Pregnancy and Postpartum'
    'Public Health, Environmental and Occupational Health'
    'Global Maternal and Child Health'
    'Pediatrics, Perinatology and Child Health'
    'Health disparities and outcomes'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Remove the now-unused rows 46-62 (old data extended further)
$lastRow = 62
$newLastRow = $values.Length + 1
if ($lastRow -gt $newLastRow) {
    $deleteRange = $ws.Range($ws.Cells.Item($newLastRow + 1, 1), $ws.Cells.Item($lastRow, 1))
    $deleteRange.EntireRow.Delete()
}

$ws.UsedRange | Out-Null
